$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.197.71'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '1.879.83'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''313.18'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = '''1.003'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '''0.5130'
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +2.55%  '
$ws.Range('D9').Value = '''0.08372'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = '''41.42'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '''6.234'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '''20.73'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').Value = '1.876.82'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '''0.00001107'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').Value = '''91.37'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').Value = '''0.06646'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '''17.75'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('D23').Value = '28.218.17'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').Value = '''11.22'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').Value = '''2.264'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '2.091.24'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = '''2.511'
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('D28').Value = '''158.50'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('D30').Value = '''125.50'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').Value = '''0.1066'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').Value = '''5.887'
$ws.Range('E33').Value = '  +5.24%  '
$ws.Range('D34').Value = '''3.589'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = '''9.744'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').Value = '''0.02457'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').Value = '''0.06556'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').Value = '''0.2194'
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('D39').Value = '''1.211'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '''0.6514'
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('D41').Value = '''5.044'
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('D42').Value = '''1.232'
$ws.Range('E42').Value = '  -0.72%  '
$ws.Range('D43').Value = '''11.31'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '''0.6114'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '''13.10'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('D47').Value = '''3.675'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').Value = '''2.019'
$ws.Range('E48').Value = '  +0.58%  '
$ws.Range('D49').Value = '''1.231'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').Value = '''121.70'
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').Value = '''78.23'
$ws.Range('E51').Value = '  -2.43%  '
